$wb = $excel.ActiveWorkbook

$newFile = "b09fff94-4a04-4c09-88e8-76b0ea303c8b.md"
$newXlfBase = "b09fff94-4a04-4c09-88e8-76b0ea303c8b.c38d42d5ea8ade9cef9101e29fdaee151d7d478e"
$oldFile = "ff65b339-b1e1-4f1d-9698-14fa462ca63d.md"
$oldXlfBase = "ff65b339-b1e1-4f1d-9698-14fa462ca63d.0a09562e4e71c938ab98e9d5225a6d0541d26152"

$newMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/932376d4014fb5dcbfbef45d0b162fdf7d62d5fd/e2e/$newFile"
$oldMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/932376d4014fb5dcbfbef45d0b162fdf7d62d5fd/e2e/$oldFile"

$newZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a9d7c3bfe8ae369fc2882def8186bd52f96707d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfBase.zh-cn.xlf"
$oldZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a9d7c3bfe8ae369fc2882def8186bd52f96707d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldXlfBase.zh-cn.xlf"

$newDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2183816714562800ffb36778a9346c551645792b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfBase.de-de.xlf"
$oldDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2183816714562800ffb36778a9346c551645792b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldXlfBase.de-de.xlf"

# ---------- Overview sheet ----------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Rows.Item(2).Insert()

$ws1.Range("A2").Value = $newFile
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-38-13 08:38:46"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $newMdUrl, "", "", $newFile)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $oldMdUrl, "", "", $oldFile)

# ---------- zh-cn sheet ----------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Rows.Item(2).Insert()

$ws2.Range("A2").Value = $newFile
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "$newXlfBase.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-13 08:38:43"
$ws2.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("I2").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $newMdUrl, "", "", $newFile)
$ws2.Hyperlinks.Add($ws2.Range("B2"), $newMdUrl, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), $newZhUrl, "", "", "$newXlfBase.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $oldMdUrl, "", "", $oldFile)
$ws2.Hyperlinks.Add($ws2.Range("B3"), $oldMdUrl, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), $oldZhUrl, "", "", "$oldXlfBase.zh-cn.xlf")

# ---------- de-de sheet ----------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Rows.Item(2).Insert()

$ws3.Range("A2").Value = $newFile
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "$newXlfBase.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-13 08:38:46"
$ws3.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("I2").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $newMdUrl, "", "", $newFile)
$ws3.Hyperlinks.Add($ws3.Range("B2"), $newMdUrl, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), $newDeUrl, "", "", "$newXlfBase.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $oldMdUrl, "", "", $oldFile)
$ws3.Hyperlinks.Add($ws3.Range("B3"), $oldMdUrl, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), $oldDeUrl, "", "", "$oldXlfBase.de-de.xlf")

Write-Output "done"
